$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (pushes Late/heading/Outstanding
# columns one to the right), matching the RBI variable-instalments layout.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Switch focus to the Repayment schedule tab and leave the selection where
# the author left it.
[void]$ws.Activate()
[void]$ws.Range("K15").Select()
